$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up a few existing text values (dash -> colon, missing space) ---
$ws.Range("B2").Value = "មេរៀនម៉ូឌុលទី ១៖ ការណែនាំអំពីគណនេយ្យភាពសង្គម (ISAF)"
$ws.Range("C2").Value = "Module 1: introduction to ISAF"
$ws.Range("C4").Value = "Module 3: Facilitating community scorecard and service provider self-assessment"

# --- Add the two new "shortcut name" columns ---
$ws.Range("D1").Value = "shortcut_name_km"
$ws.Range("E1").Value = "shortcut_name_en"
$ws.Range("D1:E1").Font.Bold = $true

$ws.Range("D2").Value = "មេរៀនម៉ូឌុលទី ១"
$ws.Range("E2").Value = "Module 1"

$ws.Range("D3").Value = "មេរៀនម៉ូឌុលទី ២"
$ws.Range("E3").Value = "Module 2"

$ws.Range("D4").Value = "មេរៀនម៉ូឌុលទី ៣"
$ws.Range("E4").Value = "Module 3"

$ws.Range("D5").Value = "មេរៀនម៉ូឌុលទី ៤"
$ws.Range("E5").Value = "Module 4"
